$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a grid of values with some cells highlighted (yellow fill).
# This edit re-shuffles which cells carry the highlight + value for a handful of
# (row, col) pairs, by swapping content+formatting between cell pairs, and
# tweaks the formatting on two more cells (H5, H16).
#
# We use an off-grid scratch cell (ZZ1) as temporary holding space so we can do
# true swaps (value + number format + fill/font) with Range.Copy, which carries
# both the value and the full cell formatting in one shot.

$scratch = $ws.Range("ZZ1")

function Swap-Cells($addrA, $addrB) {
    $a = $ws.Range($addrA)
    $b = $ws.Range($addrB)
    $a.Copy($scratch)
    $b.Copy($a)
    $scratch.Copy($b)
}

# Pairwise swaps (value + style move together)
Swap-Cells "K2" "K16"
Swap-Cells "I9" "I16"
Swap-Cells "J9" "J16"
Swap-Cells "D16" "D20"
Swap-Cells "E16" "E20"

$scratch.Clear()

# H5 loses its highlighted formatting (value stays the same) -> back to plain/default style.
$ws.Range("H5").ClearFormats()

# H16 gains a highlighted fill, with the font reset to a blank/default font (value stays
# the same). Copy the fill from a cell that already carries the yellow highlight style so
# the existing fill definition is reused, then blank out the font name.
$h16 = $ws.Range("H16")
$ws.Range("K16").Copy($h16)
$h16.Font.Name = ""
$h16.Value = -2.366257947207366
